# Score Function Loss Update
#
# The loss/prior used for the "Score Function Simulations" sheet is switched
# from a gamma distribution to a two-component Weibull mixture, and its
# parameter list is updated to match.

$wb = $excel.ActiveWorkbook

$wsScore = $wb.Worksheets.Item("Score Function Simulations")
$wsScore.Activate()

$wsScore.Range("B2").Value = "rweibullmix"
$wsScore.Range("B3").Value = "list(lambda=c(0.25,0.75), shape=c(1.1,14), scale=c(0.1,5))"

# Reflect the rows that were touched (distribution/parameters, rows 2:3) in
# the sheet's selection, landing on row 2.
$wsScore.Rows("2:3").Select() | Out-Null
